$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16:F16").NumberFormat = "@"
$ws.Cells.Item(16,3).Value() = 'Fósforo Total'
$ws.Cells.Item(16,4).Value() = ''
$ws.Cells.Item(16,5).Value() = '2.69000000'
$ws.Cells.Item(16,6).Value() = 'mg/L'

$ws.Range("C17:F17").NumberFormat = "@"
$ws.Cells.Item(17,3).Value() = 'Turbidez'
$ws.Cells.Item(17,4).Value() = ''
$ws.Cells.Item(17,5).Value() = '36.07000000'
$ws.Cells.Item(17,6).Value() = 'UNT'

$ws.Range("C18:F18").NumberFormat = "@"
$ws.Cells.Item(18,3).Value() = 'Temperatura da Água'
$ws.Cells.Item(18,4).Value() = ''
$ws.Cells.Item(18,5).Value() = '17.70000000'
$ws.Cells.Item(18,6).Value() = 'ºC'

$ws.Range("C19:F19").NumberFormat = "@"
$ws.Cells.Item(19,3).Value() = 'Escherichia coli'
$ws.Cells.Item(19,4).Value() = ''
$ws.Cells.Item(19,5).Value() = '3100000.00000000'
$ws.Cells.Item(19,6).Value() = 'UFC/100mL'

$ws.Range("C20:F20").NumberFormat = "@"
$ws.Cells.Item(20,3).Value() = 'Nitrogênio Amoniacal'
$ws.Cells.Item(20,4).Value() = ''
$ws.Cells.Item(20,5).Value() = '26.20000000'
$ws.Cells.Item(20,6).Value() = 'mg/L'

$ws.Range("C21:F21").NumberFormat = "@"
$ws.Cells.Item(21,3).Value() = 'Sólido Total'
$ws.Cells.Item(21,4).Value() = ''
$ws.Cells.Item(21,5).Value() = '398.00000000'
$ws.Cells.Item(21,6).Value() = 'mg/L'

$ws.Range("C22:F22").NumberFormat = "@"
$ws.Cells.Item(22,3).Value() = 'Carbono Orgânico Total'
$ws.Cells.Item(22,4).Value() = ''
$ws.Cells.Item(22,5).Value() = '38.80000000'
$ws.Cells.Item(22,6).Value() = 'mg/L'

$ws.Range("C23:F23").NumberFormat = "@"
$ws.Cells.Item(23,3).Value() = 'pH'
$ws.Cells.Item(23,4).Value() = ''
$ws.Cells.Item(23,5).Value() = '7.39000000'
$ws.Cells.Item(23,6).Value() = 'U.pH'

$ws.Range("C24:F24").NumberFormat = "@"
$ws.Cells.Item(24,3).Value() = 'Nitrogênio Kjeldahl'
$ws.Cells.Item(24,4).Value() = ''
$ws.Cells.Item(24,5).Value() = '31.10000000'
$ws.Cells.Item(24,6).Value() = 'mg/L'

$ws.Range("C25:F25").NumberFormat = "@"
$ws.Cells.Item(25,3).Value() = 'DQO (relativo a carbono)'
$ws.Cells.Item(25,4).Value() = ''
$ws.Cells.Item(25,5).Value() = '113.00000000'
$ws.Cells.Item(25,6).Value() = 'mg/L'

$ws.Range("C26:F26").NumberFormat = "@"
$ws.Cells.Item(26,3).Value() = 'Sólido Suspenso Total'
$ws.Cells.Item(26,4).Value() = ''
$ws.Cells.Item(26,5).Value() = '116.00000000'
$ws.Cells.Item(26,6).Value() = 'mg/L'

$ws.Range("C27:F27").NumberFormat = "@"
$ws.Cells.Item(27,3).Value() = 'DBO (5, 20)'
$ws.Cells.Item(27,4).Value() = ''
$ws.Cells.Item(27,5).Value() = '70.50000000'
$ws.Cells.Item(27,6).Value() = 'mg/L'

$ws.Range("C28:F28").NumberFormat = "@"
$ws.Cells.Item(28,3).Value() = 'Condutividade'
$ws.Cells.Item(28,4).Value() = ''
$ws.Cells.Item(28,5).Value() = '603.00000000'
$ws.Cells.Item(28,6).Value() = 'µS/cm'

$ws.Range("C29:F29").NumberFormat = "@"
$ws.Cells.Item(29,3).Value() = 'Oxigênio Dissolvido'
$ws.Cells.Item(29,4).Value() = ''
$ws.Cells.Item(29,5).Value() = '0.39000000'
$ws.Cells.Item(29,6).Value() = 'mg/L'

$ws.Range("C30:F30").NumberFormat = "@"
$ws.Cells.Item(30,3).Value() = 'Carbono Orgânico Total'
$ws.Cells.Item(30,4).Value() = ''
$ws.Cells.Item(30,5).Value() = '13.00000000'
$ws.Cells.Item(30,6).Value() = 'mg/L'

$ws.Range("C31:F31").NumberFormat = "@"
$ws.Cells.Item(31,3).Value() = 'Sólido Suspenso Total'
$ws.Cells.Item(31,4).Value() = '<'
$ws.Cells.Item(31,5).Value() = '100.00000000'
$ws.Cells.Item(31,6).Value() = 'mg/L'

$ws.Range("C32:F32").NumberFormat = "@"
$ws.Cells.Item(32,3).Value() = 'Temperatura da Água'
$ws.Cells.Item(32,4).Value() = ''
$ws.Cells.Item(32,5).Value() = '19.00000000'
$ws.Cells.Item(32,6).Value() = 'ºC'

$ws.Range("C33:F33").NumberFormat = "@"
$ws.Cells.Item(33,3).Value() = 'Escherichia coli'
$ws.Cells.Item(33,4).Value() = ''
$ws.Cells.Item(33,5).Value() = '1400000.00000000'
$ws.Cells.Item(33,6).Value() = 'UFC/100mL'

$ws.Range("C34:F34").NumberFormat = "@"
$ws.Cells.Item(34,3).Value() = 'Turbidez'
$ws.Cells.Item(34,4).Value() = ''
$ws.Cells.Item(34,5).Value() = '20.54000000'
$ws.Cells.Item(34,6).Value() = 'UNT'

$ws.Range("C35:F35").NumberFormat = "@"
$ws.Cells.Item(35,3).Value() = 'Oxigênio Dissolvido'
$ws.Cells.Item(35,4).Value() = ''
$ws.Cells.Item(35,5).Value() = '3.86000000'
$ws.Cells.Item(35,6).Value() = 'mg/L'

$ws.Range("C36:F36").NumberFormat = "@"
$ws.Cells.Item(36,3).Value() = 'Nitrogênio Kjeldahl'
$ws.Cells.Item(36,4).Value() = ''
$ws.Cells.Item(36,5).Value() = '14.80000000'
$ws.Cells.Item(36,6).Value() = 'mg/L'

$ws.Range("C37:F37").NumberFormat = "@"
$ws.Cells.Item(37,3).Value() = 'Sólido Total'
$ws.Cells.Item(37,4).Value() = ''
$ws.Cells.Item(37,5).Value() = '248.00000000'
$ws.Cells.Item(37,6).Value() = 'mg/L'

$ws.Range("C38:F38").NumberFormat = "@"
$ws.Cells.Item(38,3).Value() = 'DQO (relativo a carbono)'
$ws.Cells.Item(38,4).Value() = '<'
$ws.Cells.Item(38,5).Value() = '50.00000000'
$ws.Cells.Item(38,6).Value() = 'mg/L'

$ws.Range("C39:F39").NumberFormat = "@"
$ws.Cells.Item(39,3).Value() = 'Fósforo Total'
$ws.Cells.Item(39,4).Value() = ''
$ws.Cells.Item(39,5).Value() = '1.00000000'
$ws.Cells.Item(39,6).Value() = 'mg/L'

$ws.Range("C40:F40").NumberFormat = "@"
$ws.Cells.Item(40,3).Value() = 'Condutividade'
$ws.Cells.Item(40,4).Value() = ''
$ws.Cells.Item(40,5).Value() = '458.70000000'
$ws.Cells.Item(40,6).Value() = 'µS/cm'

$ws.Range("C41:F41").NumberFormat = "@"
$ws.Cells.Item(41,3).Value() = 'pH'
$ws.Cells.Item(41,4).Value() = ''
$ws.Cells.Item(41,5).Value() = '7.46000000'
$ws.Cells.Item(41,6).Value() = 'U.pH'

$ws.Range("C42:F42").NumberFormat = "@"
$ws.Cells.Item(42,3).Value() = 'DBO (5, 20)'
$ws.Cells.Item(42,4).Value() = ''
$ws.Cells.Item(42,5).Value() = '9.23000000'
$ws.Cells.Item(42,6).Value() = 'mg/L'

$ws.Range("C43:F43").NumberFormat = "@"
$ws.Cells.Item(43,3).Value() = 'Nitrogênio Amoniacal'
$ws.Cells.Item(43,4).Value() = ''
$ws.Cells.Item(43,5).Value() = '10.50000000'
$ws.Cells.Item(43,6).Value() = 'mg/L'

$ws.Range("C44:F44").NumberFormat = "@"
$ws.Cells.Item(44,3).Value() = 'Nitrogênio Amoniacal'
$ws.Cells.Item(44,4).Value() = ''
$ws.Cells.Item(44,5).Value() = '27.00000000'
$ws.Cells.Item(44,6).Value() = 'mg/L'

$ws.Range("C45:F45").NumberFormat = "@"
$ws.Cells.Item(45,3).Value() = 'Nitrogênio Kjeldahl'
$ws.Cells.Item(45,4).Value() = ''
$ws.Cells.Item(45,5).Value() = '31.40000000'
$ws.Cells.Item(45,6).Value() = 'mg/L'

$ws.Range("C46:F46").NumberFormat = "@"
$ws.Cells.Item(46,3).Value() = 'Carbono Orgânico Total'
$ws.Cells.Item(46,4).Value() = ''
$ws.Cells.Item(46,5).Value() = '30.00000000'
$ws.Cells.Item(46,6).Value() = 'mg/L'

$ws.Range("C47:F47").NumberFormat = "@"
$ws.Cells.Item(47,3).Value() = 'Oxigênio Dissolvido'
$ws.Cells.Item(47,4).Value() = ''
$ws.Cells.Item(47,5).Value() = '0.22000000'
$ws.Cells.Item(47,6).Value() = 'mg/L'

$ws.Range("C48:F48").NumberFormat = "@"
$ws.Cells.Item(48,3).Value() = 'Sólido Suspenso Total'
$ws.Cells.Item(48,4).Value() = '<'
$ws.Cells.Item(48,5).Value() = '100.00000000'
$ws.Cells.Item(48,6).Value() = 'mg/L'

$ws.Range("C49:F49").NumberFormat = "@"
$ws.Cells.Item(49,3).Value() = 'Temperatura da Água'
$ws.Cells.Item(49,4).Value() = ''
$ws.Cells.Item(49,5).Value() = '21.10000000'
$ws.Cells.Item(49,6).Value() = 'ºC'

$ws.Range("C50:F50").NumberFormat = "@"
$ws.Cells.Item(50,3).Value() = 'Fósforo Total'
$ws.Cells.Item(50,4).Value() = ''
$ws.Cells.Item(50,5).Value() = '2.71000000'
$ws.Cells.Item(50,6).Value() = 'mg/L'

$ws.Range("C51:F51").NumberFormat = "@"
$ws.Cells.Item(51,3).Value() = 'Sólido Total'
$ws.Cells.Item(51,4).Value() = ''
$ws.Cells.Item(51,5).Value() = '314.00000000'
$ws.Cells.Item(51,6).Value() = 'mg/L'

$ws.Range("C52:F52").NumberFormat = "@"
$ws.Cells.Item(52,3).Value() = 'Condutividade'
$ws.Cells.Item(52,4).Value() = ''
$ws.Cells.Item(52,5).Value() = '594.00000000'
$ws.Cells.Item(52,6).Value() = 'µS/cm'

$ws.Range("C53:F53").NumberFormat = "@"
$ws.Cells.Item(53,3).Value() = 'pH'
$ws.Cells.Item(53,4).Value() = ''
$ws.Cells.Item(53,5).Value() = '7.45000000'
$ws.Cells.Item(53,6).Value() = 'U.pH'

$ws.Range("C54:F54").NumberFormat = "@"
$ws.Cells.Item(54,3).Value() = 'DBO (5, 20)'
$ws.Cells.Item(54,4).Value() = ''
$ws.Cells.Item(54,5).Value() = '60.60000000'
$ws.Cells.Item(54,6).Value() = 'mg/L'

$ws.Range("C55:F55").NumberFormat = "@"
$ws.Cells.Item(55,3).Value() = 'DQO (relativo a carbono)'
$ws.Cells.Item(55,4).Value() = ''
$ws.Cells.Item(55,5).Value() = '65.00000000'
$ws.Cells.Item(55,6).Value() = 'mg/L'

$ws.Range("C56:F56").NumberFormat = "@"
$ws.Cells.Item(56,3).Value() = 'Escherichia coli'
$ws.Cells.Item(56,4).Value() = ''
$ws.Cells.Item(56,5).Value() = '3700000.00000000'
$ws.Cells.Item(56,6).Value() = 'UFC/100mL'

$ws.Range("C57:F57").NumberFormat = "@"
$ws.Cells.Item(57,3).Value() = 'Turbidez'
$ws.Cells.Item(57,4).Value() = ''
$ws.Cells.Item(57,5).Value() = '23.70000000'
$ws.Cells.Item(57,6).Value() = 'UNT'

$ws.Range("C58:F58").NumberFormat = "@"
$ws.Cells.Item(58,3).Value() = 'DBO (5, 20)'
$ws.Cells.Item(58,4).Value() = ''
$ws.Cells.Item(58,5).Value() = '24.60000000'
$ws.Cells.Item(58,6).Value() = 'mg/L'

$ws.Range("C59:F59").NumberFormat = "@"
$ws.Cells.Item(59,3).Value() = 'Carbono Orgânico Total'
$ws.Cells.Item(59,4).Value() = ''
$ws.Cells.Item(59,5).Value() = '20.20000000'
$ws.Cells.Item(59,6).Value() = 'mg/L'

$ws.Range("C60:F60").NumberFormat = "@"
$ws.Cells.Item(60,3).Value() = 'Temperatura da Água'
$ws.Cells.Item(60,4).Value() = ''
$ws.Cells.Item(60,5).Value() = '23.40000000'
$ws.Cells.Item(60,6).Value() = 'ºC'

$ws.Range("C61:F61").NumberFormat = "@"
$ws.Cells.Item(61,3).Value() = 'DQO (relativo a carbono)'
$ws.Cells.Item(61,4).Value() = '<'
$ws.Cells.Item(61,5).Value() = '50.00000000'
$ws.Cells.Item(61,6).Value() = 'mg/L'

$ws.Range("C62:F62").NumberFormat = "@"
$ws.Cells.Item(62,3).Value() = 'Sólido Total'
$ws.Cells.Item(62,4).Value() = ''
$ws.Cells.Item(62,5).Value() = '312.00000000'
$ws.Cells.Item(62,6).Value() = 'mg/L'

$ws.Range("C63:F63").NumberFormat = "@"
$ws.Cells.Item(63,3).Value() = 'Sólido Suspenso Total'
$ws.Cells.Item(63,4).Value() = '<'
$ws.Cells.Item(63,5).Value() = '100.00000000'
$ws.Cells.Item(63,6).Value() = 'mg/L'

$ws.Range("C64:F64").NumberFormat = "@"
$ws.Cells.Item(64,3).Value() = 'Nitrogênio Amoniacal'
$ws.Cells.Item(64,4).Value() = ''
$ws.Cells.Item(64,5).Value() = '12.40000000'
$ws.Cells.Item(64,6).Value() = 'mg/L'

$ws.Range("C65:F65").NumberFormat = "@"
$ws.Cells.Item(65,3).Value() = 'pH'
$ws.Cells.Item(65,4).Value() = ''
$ws.Cells.Item(65,5).Value() = '7.46000000'
$ws.Cells.Item(65,6).Value() = 'U.pH'

$ws.Range("C66:F66").NumberFormat = "@"
$ws.Cells.Item(66,3).Value() = 'Turbidez'
$ws.Cells.Item(66,4).Value() = ''
$ws.Cells.Item(66,5).Value() = '18.90000000'
$ws.Cells.Item(66,6).Value() = 'UNT'

$ws.Range("C67:F67").NumberFormat = "@"
$ws.Cells.Item(67,3).Value() = 'Oxigênio Dissolvido'
$ws.Cells.Item(67,4).Value() = ''
$ws.Cells.Item(67,5).Value() = '3.08000000'
$ws.Cells.Item(67,6).Value() = 'mg/L'

$ws.Range("C68:F68").NumberFormat = "@"
$ws.Cells.Item(68,3).Value() = 'Escherichia coli'
$ws.Cells.Item(68,4).Value() = ''
$ws.Cells.Item(68,5).Value() = '5700000.00000000'
$ws.Cells.Item(68,6).Value() = 'UFC/100mL'

$ws.Range("C69:F69").NumberFormat = "@"
$ws.Cells.Item(69,3).Value() = 'Fósforo Total'
$ws.Cells.Item(69,4).Value() = ''
$ws.Cells.Item(69,5).Value() = '1.12000000'
$ws.Cells.Item(69,6).Value() = 'mg/L'

$ws.Range("C70:F70").NumberFormat = "@"
$ws.Cells.Item(70,3).Value() = 'Condutividade'
$ws.Cells.Item(70,4).Value() = ''
$ws.Cells.Item(70,5).Value() = '524.00000000'
$ws.Cells.Item(70,6).Value() = 'µS/cm'

$ws.Range("C71:F71").NumberFormat = "@"
$ws.Cells.Item(71,3).Value() = 'Nitrogênio Kjeldahl'
$ws.Cells.Item(71,4).Value() = ''
$ws.Cells.Item(71,5).Value() = '18.70000000'
$ws.Cells.Item(71,6).Value() = 'mg/L'

$ws.Range("C72:F72").NumberFormat = "@"
$ws.Cells.Item(72,3).Value() = 'pH'
$ws.Cells.Item(72,4).Value() = ''
$ws.Cells.Item(72,5).Value() = '7.42000000'
$ws.Cells.Item(72,6).Value() = 'U.pH'

$ws.Range("C73:F73").NumberFormat = "@"
$ws.Cells.Item(73,3).Value() = 'Sólido Suspenso Total'
$ws.Cells.Item(73,4).Value() = '<'
$ws.Cells.Item(73,5).Value() = '100.00000000'
$ws.Cells.Item(73,6).Value() = 'mg/L'

$ws.Range("C74:F74").NumberFormat = "@"
$ws.Cells.Item(74,3).Value() = 'Oxigênio Dissolvido'
$ws.Cells.Item(74,4).Value() = ''
$ws.Cells.Item(74,5).Value() = '0.83000000'
$ws.Cells.Item(74,6).Value() = 'mg/L'

$ws.Range("C75:F75").NumberFormat = "@"
$ws.Cells.Item(75,3).Value() = 'Condutividade'
$ws.Cells.Item(75,4).Value() = ''
$ws.Cells.Item(75,5).Value() = '494.80000000'
$ws.Cells.Item(75,6).Value() = 'µS/cm'

$ws.Range("C76:F76").NumberFormat = "@"
$ws.Cells.Item(76,3).Value() = 'Temperatura da Água'
$ws.Cells.Item(76,4).Value() = ''
$ws.Cells.Item(76,5).Value() = '22.90000000'
$ws.Cells.Item(76,6).Value() = 'ºC'

$ws.Range("C77:F77").NumberFormat = "@"
$ws.Cells.Item(77,3).Value() = 'Carbono Orgânico Total'
$ws.Cells.Item(77,4).Value() = ''
$ws.Cells.Item(77,5).Value() = '25.60000000'
$ws.Cells.Item(77,6).Value() = 'mg/L'

$ws.Range("C78:F78").NumberFormat = "@"
$ws.Cells.Item(78,3).Value() = 'DBO (5, 20)'
$ws.Cells.Item(78,4).Value() = ''
$ws.Cells.Item(78,5).Value() = '50.40000000'
$ws.Cells.Item(78,6).Value() = 'mg/L'

$ws.Range("C79:F79").NumberFormat = "@"
$ws.Cells.Item(79,3).Value() = 'Sólido Total'
$ws.Cells.Item(79,4).Value() = ''
$ws.Cells.Item(79,5).Value() = '260.00000000'
$ws.Cells.Item(79,6).Value() = 'mg/L'

$ws.Range("C80:F80").NumberFormat = "@"
$ws.Cells.Item(80,3).Value() = 'Nitrogênio Amoniacal'
$ws.Cells.Item(80,4).Value() = ''
$ws.Cells.Item(80,5).Value() = '15.80000000'
$ws.Cells.Item(80,6).Value() = 'mg/L'

$ws.Range("C81:F81").NumberFormat = "@"
$ws.Cells.Item(81,3).Value() = 'DQO (relativo a carbono)'
$ws.Cells.Item(81,4).Value() = ''
$ws.Cells.Item(81,5).Value() = '50.40000000'
$ws.Cells.Item(81,6).Value() = 'mg/L'

$ws.Range("C82:F82").NumberFormat = "@"
$ws.Cells.Item(82,3).Value() = 'Turbidez'
$ws.Cells.Item(82,4).Value() = ''
$ws.Cells.Item(82,5).Value() = '27.12000000'
$ws.Cells.Item(82,6).Value() = 'UNT'

$ws.Range("C83:F83").NumberFormat = "@"
$ws.Cells.Item(83,3).Value() = 'Escherichia coli'
$ws.Cells.Item(83,4).Value() = ''
$ws.Cells.Item(83,5).Value() = '4100000.00000000'
$ws.Cells.Item(83,6).Value() = 'UFC/100mL'

$ws.Range("C84:F84").NumberFormat = "@"
$ws.Cells.Item(84,3).Value() = 'Fósforo Total'
$ws.Cells.Item(84,4).Value() = ''
$ws.Cells.Item(84,5).Value() = '2.00000000'
$ws.Cells.Item(84,6).Value() = 'mg/L'

$ws.Range("A85").NumberFormat = "@"
$ws.Range("B85").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C85:F85").NumberFormat = "@"
$ws.Cells.Item(85,1).Value() = 'REIM02800'
$ws.Cells.Item(85,2).Value() = 44250.38888888889
$ws.Cells.Item(85,3).Value() = 'Nitrogênio Kjeldahl'
$ws.Cells.Item(85,4).Value() = ''
$ws.Cells.Item(85,5).Value() = '22.90000000'
$ws.Cells.Item(85,6).Value() = 'mg/L'

$ws.Range("A86").NumberFormat = "@"
$ws.Range("B86").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C86:F86").NumberFormat = "@"
$ws.Cells.Item(86,1).Value() = 'REIM02800'
$ws.Cells.Item(86,2).Value() = 44313.38194444445
$ws.Cells.Item(86,3).Value() = 'Sólido Suspenso Total'
$ws.Cells.Item(86,4).Value() = '<'
$ws.Cells.Item(86,5).Value() = '100.00000000'
$ws.Cells.Item(86,6).Value() = 'mg/L'

$ws.Range("A87").NumberFormat = "@"
$ws.Range("B87").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C87:F87").NumberFormat = "@"
$ws.Cells.Item(87,1).Value() = 'REIM02800'
$ws.Cells.Item(87,2).Value() = 44313.38194444445
$ws.Cells.Item(87,3).Value() = 'pH'
$ws.Cells.Item(87,4).Value() = ''
$ws.Cells.Item(87,5).Value() = '7.41000000'
$ws.Cells.Item(87,6).Value() = 'U.pH'

$ws.Range("A88").NumberFormat = "@"
$ws.Range("B88").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C88:F88").NumberFormat = "@"
$ws.Cells.Item(88,1).Value() = 'REIM02800'
$ws.Cells.Item(88,2).Value() = 44313.38194444445
$ws.Cells.Item(88,3).Value() = 'Sólido Total'
$ws.Cells.Item(88,4).Value() = ''
$ws.Cells.Item(88,5).Value() = '302.00000000'
$ws.Cells.Item(88,6).Value() = 'mg/L'

$ws.Range("A89").NumberFormat = "@"
$ws.Range("B89").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C89:F89").NumberFormat = "@"
$ws.Cells.Item(89,1).Value() = 'REIM02800'
$ws.Cells.Item(89,2).Value() = 44313.38194444445
$ws.Cells.Item(89,3).Value() = 'Temperatura da Água'
$ws.Cells.Item(89,4).Value() = ''
$ws.Cells.Item(89,5).Value() = '20.60000000'
$ws.Cells.Item(89,6).Value() = 'ºC'

$ws.Range("A90").NumberFormat = "@"
$ws.Range("B90").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C90:F90").NumberFormat = "@"
$ws.Cells.Item(90,1).Value() = 'REIM02800'
$ws.Cells.Item(90,2).Value() = 44313.38194444445
$ws.Cells.Item(90,3).Value() = 'Turbidez'
$ws.Cells.Item(90,4).Value() = ''
$ws.Cells.Item(90,5).Value() = '64.14000000'
$ws.Cells.Item(90,6).Value() = 'UNT'

$ws.Range("A91").NumberFormat = "@"
$ws.Range("B91").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C91:F91").NumberFormat = "@"
$ws.Cells.Item(91,1).Value() = 'REIM02800'
$ws.Cells.Item(91,2).Value() = 44313.38194444445
$ws.Cells.Item(91,3).Value() = 'Nitrogênio Kjeldahl'
$ws.Cells.Item(91,4).Value() = ''
$ws.Cells.Item(91,5).Value() = '29.10000000'
$ws.Cells.Item(91,6).Value() = 'mg/L'

$ws.Range("A92").NumberFormat = "@"
$ws.Range("B92").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C92:F92").NumberFormat = "@"
$ws.Cells.Item(92,1).Value() = 'REIM02800'
$ws.Cells.Item(92,2).Value() = 44313.38194444445
$ws.Cells.Item(92,3).Value() = 'Fósforo Total'
$ws.Cells.Item(92,4).Value() = ''
$ws.Cells.Item(92,5).Value() = '3.00000000'
$ws.Cells.Item(92,6).Value() = 'mg/L'

$ws.Range("A93").NumberFormat = "@"
$ws.Range("B93").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C93:F93").NumberFormat = "@"
$ws.Cells.Item(93,1).Value() = 'REIM02800'
$ws.Cells.Item(93,2).Value() = 44313.38194444445
$ws.Cells.Item(93,3).Value() = 'Oxigênio Dissolvido'
$ws.Cells.Item(93,4).Value() = ''
$ws.Cells.Item(93,5).Value() = '0.29000000'
$ws.Cells.Item(93,6).Value() = 'mg/L'

$ws.Range("A94").NumberFormat = "@"
$ws.Range("B94").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C94:F94").NumberFormat = "@"
$ws.Cells.Item(94,1).Value() = 'REIM02800'
$ws.Cells.Item(94,2).Value() = 44313.38194444445
$ws.Cells.Item(94,3).Value() = 'Carbono Orgânico Total'
$ws.Cells.Item(94,4).Value() = ''
$ws.Cells.Item(94,5).Value() = '30.40000000'
$ws.Cells.Item(94,6).Value() = 'mg/L'

$ws.Range("A95").NumberFormat = "@"
$ws.Range("B95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C95:F95").NumberFormat = "@"
$ws.Cells.Item(95,1).Value() = 'REIM02800'
$ws.Cells.Item(95,2).Value() = 44313.38194444445
$ws.Cells.Item(95,3).Value() = 'DQO (relativo a carbono)'
$ws.Cells.Item(95,4).Value() = ''
$ws.Cells.Item(95,5).Value() = '70.20000000'
$ws.Cells.Item(95,6).Value() = 'mg/L'

$ws.Range("A96").NumberFormat = "@"
$ws.Range("B96").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C96:F96").NumberFormat = "@"
$ws.Cells.Item(96,1).Value() = 'REIM02800'
$ws.Cells.Item(96,2).Value() = 44313.38194444445
$ws.Cells.Item(96,3).Value() = 'Escherichia coli'
$ws.Cells.Item(96,4).Value() = ''
$ws.Cells.Item(96,5).Value() = '1400000.00000000'
$ws.Cells.Item(96,6).Value() = 'UFC/100mL'

$ws.Range("A97").NumberFormat = "@"
$ws.Range("B97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C97:F97").NumberFormat = "@"
$ws.Cells.Item(97,1).Value() = 'REIM02800'
$ws.Cells.Item(97,2).Value() = 44313.38194444445
$ws.Cells.Item(97,3).Value() = 'Condutividade'
$ws.Cells.Item(97,4).Value() = ''
$ws.Cells.Item(97,5).Value() = '554.00000000'
$ws.Cells.Item(97,6).Value() = 'µS/cm'

$ws.Range("A98").NumberFormat = "@"
$ws.Range("B98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C98:F98").NumberFormat = "@"
$ws.Cells.Item(98,1).Value() = 'REIM02800'
$ws.Cells.Item(98,2).Value() = 44313.38194444445
$ws.Cells.Item(98,3).Value() = 'Nitrogênio Amoniacal'
$ws.Cells.Item(98,4).Value() = ''
$ws.Cells.Item(98,5).Value() = '21.50000000'
$ws.Cells.Item(98,6).Value() = 'mg/L'

$ws.Range("A99").NumberFormat = "@"
$ws.Range("B99").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C99:F99").NumberFormat = "@"
$ws.Cells.Item(99,1).Value() = 'REIM02800'
$ws.Cells.Item(99,2).Value() = 44313.38194444445
$ws.Cells.Item(99,3).Value() = 'DBO (5, 20)'
$ws.Cells.Item(99,4).Value() = ''
$ws.Cells.Item(99,5).Value() = '68.70000000'
$ws.Cells.Item(99,6).Value() = 'mg/L'

$ws.Range("A100").NumberFormat = "@"
$ws.Range("B100").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C100:F100").NumberFormat = "@"
$ws.Cells.Item(100,1).Value() = 'REIM02800'
$ws.Cells.Item(100,2).Value() = 44364.375
$ws.Cells.Item(100,3).Value() = 'DBO (5, 20)'
$ws.Cells.Item(100,4).Value() = ''
$ws.Cells.Item(100,5).Value() = '73.40000000'
$ws.Cells.Item(100,6).Value() = 'mg/L'

$ws.Range("A101").NumberFormat = "@"
$ws.Range("B101").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C101:F101").NumberFormat = "@"
$ws.Cells.Item(101,1).Value() = 'REIM02800'
$ws.Cells.Item(101,2).Value() = 44364.375
$ws.Cells.Item(101,3).Value() = 'Escherichia coli'
$ws.Cells.Item(101,4).Value() = ''
$ws.Cells.Item(101,5).Value() = '1400000.00000000'
$ws.Cells.Item(101,6).Value() = 'UFC/100mL'

$ws.Range("A102").NumberFormat = "@"
$ws.Range("B102").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C102:F102").NumberFormat = "@"
$ws.Cells.Item(102,1).Value() = 'REIM02800'
$ws.Cells.Item(102,2).Value() = 44364.375
$ws.Cells.Item(102,3).Value() = 'Turbidez'
$ws.Cells.Item(102,4).Value() = ''
$ws.Cells.Item(102,5).Value() = '37.90000000'
$ws.Cells.Item(102,6).Value() = 'UNT'

$ws.Range("A103").NumberFormat = "@"
$ws.Range("B103").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C103:F103").NumberFormat = "@"
$ws.Cells.Item(103,1).Value() = 'REIM02800'
$ws.Cells.Item(103,2).Value() = 44364.375
$ws.Cells.Item(103,3).Value() = 'Sólido Total'
$ws.Cells.Item(103,4).Value() = ''
$ws.Cells.Item(103,5).Value() = '320.00000000'
$ws.Cells.Item(103,6).Value() = 'mg/L'

$ws.Range("A104").NumberFormat = "@"
$ws.Range("B104").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C104:F104").NumberFormat = "@"
$ws.Cells.Item(104,1).Value() = 'REIM02800'
$ws.Cells.Item(104,2).Value() = 44364.375
$ws.Cells.Item(104,3).Value() = 'Carbono Orgânico Total'
$ws.Cells.Item(104,4).Value() = ''
$ws.Cells.Item(104,5).Value() = '30.80000000'
$ws.Cells.Item(104,6).Value() = 'mg/L'

$ws.Range("A105").NumberFormat = "@"
$ws.Range("B105").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C105:F105").NumberFormat = "@"
$ws.Cells.Item(105,1).Value() = 'REIM02800'
$ws.Cells.Item(105,2).Value() = 44364.375
$ws.Cells.Item(105,3).Value() = 'Temperatura da Água'
$ws.Cells.Item(105,4).Value() = ''
$ws.Cells.Item(105,5).Value() = '17.50000000'
$ws.Cells.Item(105,6).Value() = 'ºC'

$ws.Range("A106").NumberFormat = "@"
$ws.Range("B106").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C106:F106").NumberFormat = "@"
$ws.Cells.Item(106,1).Value() = 'REIM02800'
$ws.Cells.Item(106,2).Value() = 44364.375
$ws.Cells.Item(106,3).Value() = 'Fósforo Total'
$ws.Cells.Item(106,4).Value() = ''
$ws.Cells.Item(106,5).Value() = '3.00000000'
$ws.Cells.Item(106,6).Value() = 'mg/L'

$ws.Range("A107").NumberFormat = "@"
$ws.Range("B107").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C107:F107").NumberFormat = "@"
$ws.Cells.Item(107,1).Value() = 'REIM02800'
$ws.Cells.Item(107,2).Value() = 44364.375
$ws.Cells.Item(107,3).Value() = 'Nitrogênio Kjeldahl'
$ws.Cells.Item(107,4).Value() = ''
$ws.Cells.Item(107,5).Value() = '32.10000000'
$ws.Cells.Item(107,6).Value() = 'mg/L'

$ws.Range("A108").NumberFormat = "@"
$ws.Range("B108").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C108:F108").NumberFormat = "@"
$ws.Cells.Item(108,1).Value() = 'REIM02800'
$ws.Cells.Item(108,2).Value() = 44364.375
$ws.Cells.Item(108,3).Value() = 'DQO (relativo a carbono)'
$ws.Cells.Item(108,4).Value() = ''
$ws.Cells.Item(108,5).Value() = '110.00000000'
$ws.Cells.Item(108,6).Value() = 'mg/L'

$ws.Range("A109").NumberFormat = "@"
$ws.Range("B109").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C109:F109").NumberFormat = "@"
$ws.Cells.Item(109,1).Value() = 'REIM02800'
$ws.Cells.Item(109,2).Value() = 44364.375
$ws.Cells.Item(109,3).Value() = 'Sólido Suspenso Total'
$ws.Cells.Item(109,4).Value() = '<'
$ws.Cells.Item(109,5).Value() = '100.00000000'
$ws.Cells.Item(109,6).Value() = 'mg/L'

$ws.Range("A110").NumberFormat = "@"
$ws.Range("B110").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C110:F110").NumberFormat = "@"
$ws.Cells.Item(110,1).Value() = 'REIM02800'
$ws.Cells.Item(110,2).Value() = 44364.375
$ws.Cells.Item(110,3).Value() = 'pH'
$ws.Cells.Item(110,4).Value() = ''
$ws.Cells.Item(110,5).Value() = '7.52000000'
$ws.Cells.Item(110,6).Value() = 'U.pH'

$ws.Range("A111").NumberFormat = "@"
$ws.Range("B111").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C111:F111").NumberFormat = "@"
$ws.Cells.Item(111,1).Value() = 'REIM02800'
$ws.Cells.Item(111,2).Value() = 44364.375
$ws.Cells.Item(111,3).Value() = 'Oxigênio Dissolvido'
$ws.Cells.Item(111,4).Value() = ''
$ws.Cells.Item(111,5).Value() = '0.41000000'
$ws.Cells.Item(111,6).Value() = 'mg/L'

$ws.Range("A112").NumberFormat = "@"
$ws.Range("B112").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C112:F112").NumberFormat = "@"
$ws.Cells.Item(112,1).Value() = 'REIM02800'
$ws.Cells.Item(112,2).Value() = 44364.375
$ws.Cells.Item(112,3).Value() = 'Nitrogênio Amoniacal'
$ws.Cells.Item(112,4).Value() = ''
$ws.Cells.Item(112,5).Value() = '21.90000000'
$ws.Cells.Item(112,6).Value() = 'mg/L'

$ws.Range("A113").NumberFormat = "@"
$ws.Range("B113").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C113:F113").NumberFormat = "@"
$ws.Cells.Item(113,1).Value() = 'REIM02800'
$ws.Cells.Item(113,2).Value() = 44364.375
$ws.Cells.Item(113,3).Value() = 'Condutividade'
$ws.Cells.Item(113,4).Value() = ''
$ws.Cells.Item(113,5).Value() = '547.00000000'
$ws.Cells.Item(113,6).Value() = 'µS/cm'

$ws.Range("A114").NumberFormat = "@"
$ws.Range("B114").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C114:F114").NumberFormat = "@"
$ws.Cells.Item(114,1).Value() = 'REIM02800'
$ws.Cells.Item(114,2).Value() = 44426.38541666666
$ws.Cells.Item(114,3).Value() = 'Temperatura da Água'
$ws.Cells.Item(114,4).Value() = ''
$ws.Cells.Item(114,5).Value() = '17.90000000'
$ws.Cells.Item(114,6).Value() = 'ºC'

$ws.Range("A115").NumberFormat = "@"
$ws.Range("B115").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C115:F115").NumberFormat = "@"
$ws.Cells.Item(115,1).Value() = 'REIM02800'
$ws.Cells.Item(115,2).Value() = 44426.38541666666
$ws.Cells.Item(115,3).Value() = 'Oxigênio Dissolvido'
$ws.Cells.Item(115,4).Value() = ''
$ws.Cells.Item(115,5).Value() = '0.53000000'
$ws.Cells.Item(115,6).Value() = 'mg/L'

$ws.Range("A116").NumberFormat = "@"
$ws.Range("B116").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C116:F116").NumberFormat = "@"
$ws.Cells.Item(116,1).Value() = 'REIM02800'
$ws.Cells.Item(116,2).Value() = 44426.38541666666
$ws.Cells.Item(116,3).Value() = 'DQO (relativo a carbono)'
$ws.Cells.Item(116,4).Value() = ''
$ws.Cells.Item(116,5).Value() = '98.50000000'
$ws.Cells.Item(116,6).Value() = 'mg/L'

$ws.Range("A117").NumberFormat = "@"
$ws.Range("B117").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C117:F117").NumberFormat = "@"
$ws.Cells.Item(117,1).Value() = 'REIM02800'
$ws.Cells.Item(117,2).Value() = 44426.38541666666
$ws.Cells.Item(117,3).Value() = 'DBO (5, 20)'
$ws.Cells.Item(117,4).Value() = ''
$ws.Cells.Item(117,5).Value() = '69.50000000'
$ws.Cells.Item(117,6).Value() = 'mg/L'

$ws.Range("A118").NumberFormat = "@"
$ws.Range("B118").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C118:F118").NumberFormat = "@"
$ws.Cells.Item(118,1).Value() = 'REIM02800'
$ws.Cells.Item(118,2).Value() = 44426.38541666666
$ws.Cells.Item(118,3).Value() = 'Condutividade'
$ws.Cells.Item(118,4).Value() = ''
$ws.Cells.Item(118,5).Value() = '604.00000000'
$ws.Cells.Item(118,6).Value() = 'µS/cm'

$ws.Range("A119").NumberFormat = "@"
$ws.Range("B119").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C119:F119").NumberFormat = "@"
$ws.Cells.Item(119,1).Value() = 'REIM02800'
$ws.Cells.Item(119,2).Value() = 44426.38541666666
$ws.Cells.Item(119,3).Value() = 'Carbono Orgânico Total'
$ws.Cells.Item(119,4).Value() = ''
$ws.Cells.Item(119,5).Value() = '36.90000000'
$ws.Cells.Item(119,6).Value() = 'mg/L'

$ws.Range("A120").NumberFormat = "@"
$ws.Range("B120").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C120:F120").NumberFormat = "@"
$ws.Cells.Item(120,1).Value() = 'REIM02800'
$ws.Cells.Item(120,2).Value() = 44426.38541666666
$ws.Cells.Item(120,3).Value() = 'Fósforo Total'
$ws.Cells.Item(120,4).Value() = ''
$ws.Cells.Item(120,5).Value() = '2.78000000'
$ws.Cells.Item(120,6).Value() = 'mg/L'

$ws.Range("A121").NumberFormat = "@"
$ws.Range("B121").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C121:F121").NumberFormat = "@"
$ws.Cells.Item(121,1).Value() = 'REIM02800'
$ws.Cells.Item(121,2).Value() = 44426.38541666666
$ws.Cells.Item(121,3).Value() = 'Sólido Total'
$ws.Cells.Item(121,4).Value() = ''
$ws.Cells.Item(121,5).Value() = '338.00000000'
$ws.Cells.Item(121,6).Value() = 'mg/L'

$ws.Range("A122").NumberFormat = "@"
$ws.Range("B122").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C122:F122").NumberFormat = "@"
$ws.Cells.Item(122,1).Value() = 'REIM02800'
$ws.Cells.Item(122,2).Value() = 44426.38541666666
$ws.Cells.Item(122,3).Value() = 'Nitrogênio Amoniacal'
$ws.Cells.Item(122,4).Value() = ''
$ws.Cells.Item(122,5).Value() = '22.70000000'
$ws.Cells.Item(122,6).Value() = 'mg/L'

$ws.Range("A123").NumberFormat = "@"
$ws.Range("B123").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C123:F123").NumberFormat = "@"
$ws.Cells.Item(123,1).Value() = 'REIM02800'
$ws.Cells.Item(123,2).Value() = 44426.38541666666
$ws.Cells.Item(123,3).Value() = 'pH'
$ws.Cells.Item(123,4).Value() = ''
$ws.Cells.Item(123,5).Value() = '7.45000000'
$ws.Cells.Item(123,6).Value() = 'U.pH'

$ws.Range("A124").NumberFormat = "@"
$ws.Range("B124").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C124:F124").NumberFormat = "@"
$ws.Cells.Item(124,1).Value() = 'REIM02800'
$ws.Cells.Item(124,2).Value() = 44426.38541666666
$ws.Cells.Item(124,3).Value() = 'Sólido Suspenso Total'
$ws.Cells.Item(124,4).Value() = '<'
$ws.Cells.Item(124,5).Value() = '100.00000000'
$ws.Cells.Item(124,6).Value() = 'mg/L'

$ws.Range("A125").NumberFormat = "@"
$ws.Range("B125").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C125:F125").NumberFormat = "@"
$ws.Cells.Item(125,1).Value() = 'REIM02800'
$ws.Cells.Item(125,2).Value() = 44426.38541666666
$ws.Cells.Item(125,3).Value() = 'Nitrogênio Kjeldahl'
$ws.Cells.Item(125,4).Value() = ''
$ws.Cells.Item(125,5).Value() = '29.90000000'
$ws.Cells.Item(125,6).Value() = 'mg/L'

$ws.Range("A126").NumberFormat = "@"
$ws.Range("B126").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C126:F126").NumberFormat = "@"
$ws.Cells.Item(126,1).Value() = 'REIM02800'
$ws.Cells.Item(126,2).Value() = 44426.38541666666
$ws.Cells.Item(126,3).Value() = 'Escherichia coli'
$ws.Cells.Item(126,4).Value() = ''
$ws.Cells.Item(126,5).Value() = '2500000.00000000'
$ws.Cells.Item(126,6).Value() = 'UFC/100mL'

$ws.Range("A127").NumberFormat = "@"
$ws.Range("B127").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C127:F127").NumberFormat = "@"
$ws.Cells.Item(127,1).Value() = 'REIM02800'
$ws.Cells.Item(127,2).Value() = 44426.38541666666
$ws.Cells.Item(127,3).Value() = 'Turbidez'
$ws.Cells.Item(127,4).Value() = ''
$ws.Cells.Item(127,5).Value() = '31.40000000'
$ws.Cells.Item(127,6).Value() = 'UNT'

$ws.Range("A128").NumberFormat = "@"
$ws.Range("B128").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C128:F128").NumberFormat = "@"
$ws.Cells.Item(128,1).Value() = 'REIM02800'
$ws.Cells.Item(128,2).Value() = 44489.375
$ws.Cells.Item(128,3).Value() = 'Condutividade'
$ws.Cells.Item(128,4).Value() = ''
$ws.Cells.Item(128,5).Value() = '472.80000000'
$ws.Cells.Item(128,6).Value() = 'µS/cm'

$ws.Range("A129").NumberFormat = "@"
$ws.Range("B129").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C129:F129").NumberFormat = "@"
$ws.Cells.Item(129,1).Value() = 'REIM02800'
$ws.Cells.Item(129,2).Value() = 44489.375
$ws.Cells.Item(129,3).Value() = 'Turbidez'
$ws.Cells.Item(129,4).Value() = ''
$ws.Cells.Item(129,5).Value() = '13.30000000'
$ws.Cells.Item(129,6).Value() = 'UNT'

$ws.Range("A130").NumberFormat = "@"
$ws.Range("B130").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C130:F130").NumberFormat = "@"
$ws.Cells.Item(130,1).Value() = 'REIM02800'
$ws.Cells.Item(130,2).Value() = 44489.375
$ws.Cells.Item(130,3).Value() = 'Escherichia coli'
$ws.Cells.Item(130,4).Value() = ''
$ws.Cells.Item(130,5).Value() = '1200000.00000000'
$ws.Cells.Item(130,6).Value() = 'UFC/100mL'

$ws.Range("A131").NumberFormat = "@"
$ws.Range("B131").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C131:F131").NumberFormat = "@"
$ws.Cells.Item(131,1).Value() = 'REIM02800'
$ws.Cells.Item(131,2).Value() = 44489.375
$ws.Cells.Item(131,3).Value() = 'Carbono Orgânico Total'
$ws.Cells.Item(131,4).Value() = ''
$ws.Cells.Item(131,5).Value() = '15.90000000'
$ws.Cells.Item(131,6).Value() = 'mg/L'

$ws.Range("A132").NumberFormat = "@"
$ws.Range("B132").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C132:F132").NumberFormat = "@"
$ws.Cells.Item(132,1).Value() = 'REIM02800'
$ws.Cells.Item(132,2).Value() = 44489.375
$ws.Cells.Item(132,3).Value() = 'Temperatura da Água'
$ws.Cells.Item(132,4).Value() = ''
$ws.Cells.Item(132,5).Value() = '18.00000000'
$ws.Cells.Item(132,6).Value() = 'ºC'

$ws.Range("A133").NumberFormat = "@"
$ws.Range("B133").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C133:F133").NumberFormat = "@"
$ws.Cells.Item(133,1).Value() = 'REIM02800'
$ws.Cells.Item(133,2).Value() = 44489.375
$ws.Cells.Item(133,3).Value() = 'Oxigênio Dissolvido'
$ws.Cells.Item(133,4).Value() = ''
$ws.Cells.Item(133,5).Value() = '2.73000000'
$ws.Cells.Item(133,6).Value() = 'mg/L'

$ws.Range("A134").NumberFormat = "@"
$ws.Range("B134").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C134:F134").NumberFormat = "@"
$ws.Cells.Item(134,1).Value() = 'REIM02800'
$ws.Cells.Item(134,2).Value() = 44489.375
$ws.Cells.Item(134,3).Value() = 'Fósforo Total'
$ws.Cells.Item(134,4).Value() = ''
$ws.Cells.Item(134,5).Value() = '1.37000000'
$ws.Cells.Item(134,6).Value() = 'mg/L'

$ws.Range("A135").NumberFormat = "@"
$ws.Range("B135").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C135:F135").NumberFormat = "@"
$ws.Cells.Item(135,1).Value() = 'REIM02800'
$ws.Cells.Item(135,2).Value() = 44489.375
$ws.Cells.Item(135,3).Value() = 'Sólido Total'
$ws.Cells.Item(135,4).Value() = ''
$ws.Cells.Item(135,5).Value() = '248.00000000'
$ws.Cells.Item(135,6).Value() = 'mg/L'

$ws.Range("A136").NumberFormat = "@"
$ws.Range("B136").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C136:F136").NumberFormat = "@"
$ws.Cells.Item(136,1).Value() = 'REIM02800'
$ws.Cells.Item(136,2).Value() = 44489.375
$ws.Cells.Item(136,3).Value() = 'Nitrogênio Amoniacal'
$ws.Cells.Item(136,4).Value() = ''
$ws.Cells.Item(136,5).Value() = '13.20000000'
$ws.Cells.Item(136,6).Value() = 'mg/L'

$ws.Range("A137").NumberFormat = "@"
$ws.Range("B137").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C137:F137").NumberFormat = "@"
$ws.Cells.Item(137,1).Value() = 'REIM02800'
$ws.Cells.Item(137,2).Value() = 44489.375
$ws.Cells.Item(137,3).Value() = 'Sólido Suspenso Total'
$ws.Cells.Item(137,4).Value() = '<'
$ws.Cells.Item(137,5).Value() = '100.00000000'
$ws.Cells.Item(137,6).Value() = 'mg/L'

$ws.Range("A138").NumberFormat = "@"
$ws.Range("B138").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C138:F138").NumberFormat = "@"
$ws.Cells.Item(138,1).Value() = 'REIM02800'
$ws.Cells.Item(138,2).Value() = 44489.375
$ws.Cells.Item(138,3).Value() = 'Nitrogênio Kjeldahl'
$ws.Cells.Item(138,4).Value() = ''
$ws.Cells.Item(138,5).Value() = '17.80000000'
$ws.Cells.Item(138,6).Value() = 'mg/L'

$ws.Range("A139").NumberFormat = "@"
$ws.Range("B139").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C139:F139").NumberFormat = "@"
$ws.Cells.Item(139,1).Value() = 'REIM02800'
$ws.Cells.Item(139,2).Value() = 44489.375
$ws.Cells.Item(139,3).Value() = 'DBO (5, 20)'
$ws.Cells.Item(139,4).Value() = ''
$ws.Cells.Item(139,5).Value() = '20.80000000'
$ws.Cells.Item(139,6).Value() = 'mg/L'

$ws.Range("A140").NumberFormat = "@"
$ws.Range("B140").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C140:F140").NumberFormat = "@"
$ws.Cells.Item(140,1).Value() = 'REIM02800'
$ws.Cells.Item(140,2).Value() = 44489.375
$ws.Cells.Item(140,3).Value() = 'pH'
$ws.Cells.Item(140,4).Value() = ''
$ws.Cells.Item(140,5).Value() = '7.39000000'
$ws.Cells.Item(140,6).Value() = 'U.pH'

